$wb = $excel.ActiveWorkbook

# ---- PIR sheet (Worksheets.Item(2)) ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("A93:A106").NumberFormat = "@"

$ws.Cells.Item(93, 1).Value = '2026-01-28'
$ws.Cells.Item(93, 2).Value = '16:18:11'
$ws.Cells.Item(93, 3).Value = '16:00'
$ws.Cells.Item(93, 4).Value = 'Bathroom'
$ws.Cells.Item(93, 5).Value = 'No Motion'
$ws.Cells.Item(93, 6).Value = 'Inactive'
$ws.Cells.Item(94, 1).Value = '2026-01-28'
$ws.Cells.Item(94, 2).Value = '16:18:12'
$ws.Cells.Item(94, 3).Value = '16:00'
$ws.Cells.Item(94, 4).Value = 'Bathroom'
$ws.Cells.Item(94, 5).Value = 'No Motion'
$ws.Cells.Item(94, 6).Value = 'Inactive'
$ws.Cells.Item(95, 1).Value = '2026-01-28'
$ws.Cells.Item(95, 2).Value = '16:18:15'
$ws.Cells.Item(95, 3).Value = '16:00'
$ws.Cells.Item(95, 4).Value = 'Bathroom'
$ws.Cells.Item(95, 5).Value = 'No Motion'
$ws.Cells.Item(95, 6).Value = 'Inactive'
$ws.Cells.Item(96, 1).Value = '2026-01-28'
$ws.Cells.Item(96, 2).Value = '16:18:20'
$ws.Cells.Item(96, 3).Value = '16:00'
$ws.Cells.Item(96, 4).Value = 'Bathroom'
$ws.Cells.Item(96, 5).Value = 'No Motion'
$ws.Cells.Item(96, 6).Value = 'Inactive'
$ws.Cells.Item(97, 1).Value = '2026-01-28'
$ws.Cells.Item(97, 2).Value = '16:18:25'
$ws.Cells.Item(97, 3).Value = '16:00'
$ws.Cells.Item(97, 4).Value = 'Bathroom'
$ws.Cells.Item(97, 5).Value = 'No Motion'
$ws.Cells.Item(97, 6).Value = 'Inactive'
$ws.Cells.Item(98, 1).Value = '2026-01-28'
$ws.Cells.Item(98, 2).Value = '16:18:30'
$ws.Cells.Item(98, 3).Value = '16:00'
$ws.Cells.Item(98, 4).Value = 'Bathroom'
$ws.Cells.Item(98, 5).Value = 'No Motion'
$ws.Cells.Item(98, 6).Value = 'Inactive'
$ws.Cells.Item(99, 1).Value = '2026-01-28'
$ws.Cells.Item(99, 2).Value = '16:18:35'
$ws.Cells.Item(99, 3).Value = '16:00'
$ws.Cells.Item(99, 4).Value = 'Bathroom'
$ws.Cells.Item(99, 5).Value = 'No Motion'
$ws.Cells.Item(99, 6).Value = 'Inactive'
$ws.Cells.Item(100, 1).Value = '2026-01-28'
$ws.Cells.Item(100, 2).Value = '16:18:40'
$ws.Cells.Item(100, 3).Value = '16:00'
$ws.Cells.Item(100, 4).Value = 'Bathroom'
$ws.Cells.Item(100, 5).Value = 'No Motion'
$ws.Cells.Item(100, 6).Value = 'Inactive'
$ws.Cells.Item(101, 1).Value = '2026-01-28'
$ws.Cells.Item(101, 2).Value = '16:18:45'
$ws.Cells.Item(101, 3).Value = '16:00'
$ws.Cells.Item(101, 4).Value = 'Bathroom'
$ws.Cells.Item(101, 5).Value = 'No Motion'
$ws.Cells.Item(101, 6).Value = 'Inactive'
$ws.Cells.Item(102, 1).Value = '2026-01-28'
$ws.Cells.Item(102, 2).Value = '16:18:50'
$ws.Cells.Item(102, 3).Value = '16:00'
$ws.Cells.Item(102, 4).Value = 'Bathroom'
$ws.Cells.Item(102, 5).Value = 'No Motion'
$ws.Cells.Item(102, 6).Value = 'Inactive'
$ws.Cells.Item(103, 1).Value = '2026-01-28'
$ws.Cells.Item(103, 2).Value = '16:18:55'
$ws.Cells.Item(103, 3).Value = '16:00'
$ws.Cells.Item(103, 4).Value = 'Bathroom'
$ws.Cells.Item(103, 5).Value = 'No Motion'
$ws.Cells.Item(103, 6).Value = 'Inactive'
$ws.Cells.Item(104, 1).Value = '2026-01-28'
$ws.Cells.Item(104, 2).Value = '16:19:00'
$ws.Cells.Item(104, 3).Value = '16:00'
$ws.Cells.Item(104, 4).Value = 'Bathroom'
$ws.Cells.Item(104, 5).Value = 'No Motion'
$ws.Cells.Item(104, 6).Value = 'Inactive'
$ws.Cells.Item(105, 1).Value = '2026-01-28'
$ws.Cells.Item(105, 2).Value = '16:19:05'
$ws.Cells.Item(105, 3).Value = '16:00'
$ws.Cells.Item(105, 4).Value = 'Bathroom'
$ws.Cells.Item(105, 5).Value = 'No Motion'
$ws.Cells.Item(105, 6).Value = 'Inactive'
$ws.Cells.Item(106, 1).Value = '2026-01-28'
$ws.Cells.Item(106, 2).Value = '16:19:10'
$ws.Cells.Item(106, 3).Value = '16:00'
$ws.Cells.Item(106, 4).Value = 'Bathroom'
$ws.Cells.Item(106, 5).Value = 'No Motion'
$ws.Cells.Item(106, 6).Value = 'Inactive'

# ---- Humidity sheet (Worksheets.Item(3)) ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("A94:A107").NumberFormat = "@"
$ws.Range("E94:E107").NumberFormat = "@"

$ws.Cells.Item(94, 1).Value = '2026-01-28'
$ws.Cells.Item(94, 2).Value = '16:18:11'
$ws.Cells.Item(94, 3).Value = '16:00'
$ws.Cells.Item(94, 4).Value = 'Bathroom'
$ws.Cells.Item(94, 5).Value = '87.5%'
$ws.Cells.Item(94, 6).Value = 'Active'
$ws.Cells.Item(95, 1).Value = '2026-01-28'
$ws.Cells.Item(95, 2).Value = '16:18:12'
$ws.Cells.Item(95, 3).Value = '16:00'
$ws.Cells.Item(95, 4).Value = 'Bathroom'
$ws.Cells.Item(95, 5).Value = '88.4%'
$ws.Cells.Item(95, 6).Value = 'Active'
$ws.Cells.Item(96, 1).Value = '2026-01-28'
$ws.Cells.Item(96, 2).Value = '16:18:15'
$ws.Cells.Item(96, 3).Value = '16:00'
$ws.Cells.Item(96, 4).Value = 'Bathroom'
$ws.Cells.Item(96, 5).Value = '87.5%'
$ws.Cells.Item(96, 6).Value = 'Active'
$ws.Cells.Item(97, 1).Value = '2026-01-28'
$ws.Cells.Item(97, 2).Value = '16:18:19'
$ws.Cells.Item(97, 3).Value = '16:00'
$ws.Cells.Item(97, 4).Value = 'Bathroom'
$ws.Cells.Item(97, 5).Value = '88.5%'
$ws.Cells.Item(97, 6).Value = 'Active'
$ws.Cells.Item(98, 1).Value = '2026-01-28'
$ws.Cells.Item(98, 2).Value = '16:18:27'
$ws.Cells.Item(98, 3).Value = '16:00'
$ws.Cells.Item(98, 4).Value = 'Bathroom'
$ws.Cells.Item(98, 5).Value = '88.4%'
$ws.Cells.Item(98, 6).Value = 'Active'
$ws.Cells.Item(99, 1).Value = '2026-01-28'
$ws.Cells.Item(99, 2).Value = '16:18:31'
$ws.Cells.Item(99, 3).Value = '16:00'
$ws.Cells.Item(99, 4).Value = 'Bathroom'
$ws.Cells.Item(99, 5).Value = '88.5%'
$ws.Cells.Item(99, 6).Value = 'Active'
$ws.Cells.Item(100, 1).Value = '2026-01-28'
$ws.Cells.Item(100, 2).Value = '16:18:35'
$ws.Cells.Item(100, 3).Value = '16:00'
$ws.Cells.Item(100, 4).Value = 'Bathroom'
$ws.Cells.Item(100, 5).Value = '87.5%'
$ws.Cells.Item(100, 6).Value = 'Active'
$ws.Cells.Item(101, 1).Value = '2026-01-28'
$ws.Cells.Item(101, 2).Value = '16:18:39'
$ws.Cells.Item(101, 3).Value = '16:00'
$ws.Cells.Item(101, 4).Value = 'Bathroom'
$ws.Cells.Item(101, 5).Value = '88.5%'
$ws.Cells.Item(101, 6).Value = 'Active'
$ws.Cells.Item(102, 1).Value = '2026-01-28'
$ws.Cells.Item(102, 2).Value = '16:18:43'
$ws.Cells.Item(102, 3).Value = '16:00'
$ws.Cells.Item(102, 4).Value = 'Bathroom'
$ws.Cells.Item(102, 5).Value = '88.5%'
$ws.Cells.Item(102, 6).Value = 'Active'
$ws.Cells.Item(103, 1).Value = '2026-01-28'
$ws.Cells.Item(103, 2).Value = '16:18:47'
$ws.Cells.Item(103, 3).Value = '16:00'
$ws.Cells.Item(103, 4).Value = 'Bathroom'
$ws.Cells.Item(103, 5).Value = '87.5%'
$ws.Cells.Item(103, 6).Value = 'Active'
$ws.Cells.Item(104, 1).Value = '2026-01-28'
$ws.Cells.Item(104, 2).Value = '16:18:52'
$ws.Cells.Item(104, 3).Value = '16:00'
$ws.Cells.Item(104, 4).Value = 'Bathroom'
$ws.Cells.Item(104, 5).Value = '88.5%'
$ws.Cells.Item(104, 6).Value = 'Active'
$ws.Cells.Item(105, 1).Value = '2026-01-28'
$ws.Cells.Item(105, 2).Value = '16:18:56'
$ws.Cells.Item(105, 3).Value = '16:00'
$ws.Cells.Item(105, 4).Value = 'Bathroom'
$ws.Cells.Item(105, 5).Value = '87.5%'
$ws.Cells.Item(105, 6).Value = 'Active'
$ws.Cells.Item(106, 1).Value = '2026-01-28'
$ws.Cells.Item(106, 2).Value = '16:19:00'
$ws.Cells.Item(106, 3).Value = '16:00'
$ws.Cells.Item(106, 4).Value = 'Bathroom'
$ws.Cells.Item(106, 5).Value = '88.5%'
$ws.Cells.Item(106, 6).Value = 'Active'
$ws.Cells.Item(107, 1).Value = '2026-01-28'
$ws.Cells.Item(107, 2).Value = '16:19:08'
$ws.Cells.Item(107, 3).Value = '16:00'
$ws.Cells.Item(107, 4).Value = 'Bathroom'
$ws.Cells.Item(107, 5).Value = '87.5%'
$ws.Cells.Item(107, 6).Value = 'Active'

# ---- Temperature sheet (Worksheets.Item(4)) ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("A94:A107").NumberFormat = "@"

$ws.Cells.Item(94, 1).Value = '2026-01-28'
$ws.Cells.Item(94, 2).Value = '16:18:12'
$ws.Cells.Item(94, 3).Value = '16:00'
$ws.Cells.Item(94, 4).Value = 'Bathroom'
$ws.Cells.Item(94, 5).Value = '22.8C'
$ws.Cells.Item(94, 6).Value = 'Active'
$ws.Cells.Item(95, 1).Value = '2026-01-28'
$ws.Cells.Item(95, 2).Value = '16:18:13'
$ws.Cells.Item(95, 3).Value = '16:00'
$ws.Cells.Item(95, 4).Value = 'Bathroom'
$ws.Cells.Item(95, 5).Value = '22.8C'
$ws.Cells.Item(95, 6).Value = 'Active'
$ws.Cells.Item(96, 1).Value = '2026-01-28'
$ws.Cells.Item(96, 2).Value = '16:18:16'
$ws.Cells.Item(96, 3).Value = '16:00'
$ws.Cells.Item(96, 4).Value = 'Bathroom'
$ws.Cells.Item(96, 5).Value = '22.8C'
$ws.Cells.Item(96, 6).Value = 'Active'
$ws.Cells.Item(97, 1).Value = '2026-01-28'
$ws.Cells.Item(97, 2).Value = '16:18:20'
$ws.Cells.Item(97, 3).Value = '16:00'
$ws.Cells.Item(97, 4).Value = 'Bathroom'
$ws.Cells.Item(97, 5).Value = '22.8C'
$ws.Cells.Item(97, 6).Value = 'Active'
$ws.Cells.Item(98, 1).Value = '2026-01-28'
$ws.Cells.Item(98, 2).Value = '16:18:28'
$ws.Cells.Item(98, 3).Value = '16:00'
$ws.Cells.Item(98, 4).Value = 'Bathroom'
$ws.Cells.Item(98, 5).Value = '22.7C'
$ws.Cells.Item(98, 6).Value = 'Active'
$ws.Cells.Item(99, 1).Value = '2026-01-28'
$ws.Cells.Item(99, 2).Value = '16:18:32'
$ws.Cells.Item(99, 3).Value = '16:00'
$ws.Cells.Item(99, 4).Value = 'Bathroom'
$ws.Cells.Item(99, 5).Value = '22.8C'
$ws.Cells.Item(99, 6).Value = 'Active'
$ws.Cells.Item(100, 1).Value = '2026-01-28'
$ws.Cells.Item(100, 2).Value = '16:18:36'
$ws.Cells.Item(100, 3).Value = '16:00'
$ws.Cells.Item(100, 4).Value = 'Bathroom'
$ws.Cells.Item(100, 5).Value = '22.8C'
$ws.Cells.Item(100, 6).Value = 'Active'
$ws.Cells.Item(101, 1).Value = '2026-01-28'
$ws.Cells.Item(101, 2).Value = '16:18:40'
$ws.Cells.Item(101, 3).Value = '16:00'
$ws.Cells.Item(101, 4).Value = 'Bathroom'
$ws.Cells.Item(101, 5).Value = '22.8C'
$ws.Cells.Item(101, 6).Value = 'Active'
$ws.Cells.Item(102, 1).Value = '2026-01-28'
$ws.Cells.Item(102, 2).Value = '16:18:44'
$ws.Cells.Item(102, 3).Value = '16:00'
$ws.Cells.Item(102, 4).Value = 'Bathroom'
$ws.Cells.Item(102, 5).Value = '22.8C'
$ws.Cells.Item(102, 6).Value = 'Active'
$ws.Cells.Item(103, 1).Value = '2026-01-28'
$ws.Cells.Item(103, 2).Value = '16:18:48'
$ws.Cells.Item(103, 3).Value = '16:00'
$ws.Cells.Item(103, 4).Value = 'Bathroom'
$ws.Cells.Item(103, 5).Value = '22.7C'
$ws.Cells.Item(103, 6).Value = 'Active'
$ws.Cells.Item(104, 1).Value = '2026-01-28'
$ws.Cells.Item(104, 2).Value = '16:18:52'
$ws.Cells.Item(104, 3).Value = '16:00'
$ws.Cells.Item(104, 4).Value = 'Bathroom'
$ws.Cells.Item(104, 5).Value = '22.8C'
$ws.Cells.Item(104, 6).Value = 'Active'
$ws.Cells.Item(105, 1).Value = '2026-01-28'
$ws.Cells.Item(105, 2).Value = '16:18:56'
$ws.Cells.Item(105, 3).Value = '16:00'
$ws.Cells.Item(105, 4).Value = 'Bathroom'
$ws.Cells.Item(105, 5).Value = '22.7C'
$ws.Cells.Item(105, 6).Value = 'Active'
$ws.Cells.Item(106, 1).Value = '2026-01-28'
$ws.Cells.Item(106, 2).Value = '16:19:00'
$ws.Cells.Item(106, 3).Value = '16:00'
$ws.Cells.Item(106, 4).Value = 'Bathroom'
$ws.Cells.Item(106, 5).Value = '22.8C'
$ws.Cells.Item(106, 6).Value = 'Active'
$ws.Cells.Item(107, 1).Value = '2026-01-28'
$ws.Cells.Item(107, 2).Value = '16:19:08'
$ws.Cells.Item(107, 3).Value = '16:00'
$ws.Cells.Item(107, 4).Value = 'Bathroom'
$ws.Cells.Item(107, 5).Value = '22.7C'
$ws.Cells.Item(107, 6).Value = 'Active'
